$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H33 value to the re-saved precision from the diff
$ws.Range("H33").Value = 45518.93175145833

# Row 34
$ws.Range("A34").Value = "Медиамонитор"
$ws.Range("B34").Value = "ПРГС.465000.024"
$ws.Range("C34").Value = 768
$ws.Range("D34").Value = 1023
$ws.Range("E34").Value = "10.8.12.0"
$ws.Range("F34").Value = "10.8.12.255"
$ws.Range("G34").Value = 256
$ws.Range("H34").Value = 45518.93761497685
$ws.Range("H34").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 35
$ws.Range("A35").Value = "Домик для Мышки Норушки"
$ws.Range("B35").Value = "с трубой"
$ws.Range("C35").Value = 769
$ws.Range("D35").Value = 776
$ws.Range("E35").Value = "10.6.23.0"
$ws.Range("F35").Value = "10.6.23.7"
$ws.Range("G35").Value = 8
$ws.Range("H35").Value = 45519.71073408565
$ws.Range("H35").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 36
$ws.Range("A36").Value = "Коммутатор"
$ws.Range("B36").Value = "ПРГС.465000.028"
$ws.Range("C36").Value = 969
$ws.Range("D36").Value = 974
$ws.Range("E36").Value = "10.8.15.0"
$ws.Range("F36").Value = "10.8.15.5"
$ws.Range("G36").Value = 6
$ws.Range("H36").Value = 45519.7180065162
$ws.Range("H36").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 37
$ws.Range("A37").Value = "Домик для Мышки Норушки"
$ws.Range("B37").Value = "без трубы"
$ws.Range("C37").Value = 2560
$ws.Range("D37").Value = 2569
$ws.Range("E37").Value = "10.6.20.1"
$ws.Range("F37").Value = "10.6.20.10"
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 45519.78871003472
$ws.Range("H37").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 38
$ws.Range("A38").Value = "Домик для Мышки Норушки"
$ws.Range("B38").Value = "без трубы"
$ws.Range("C38").Value = 2570
$ws.Range("D38").Value = 2580
$ws.Range("E38").Value = "10.6.20.11"
$ws.Range("F38").Value = "10.6.20.21"
$ws.Range("G38").Value = 11
$ws.Range("H38").Value = 45519.78922873842
$ws.Range("H38").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 39
$ws.Range("A39").Value = "Домик для Мышки Норушки"
$ws.Range("B39").Value = "без трубы"
$ws.Range("C39").Value = 2581
$ws.Range("D39").Value = 2880
$ws.Range("E39").Value = "10.6.20.22"
$ws.Range("F39").Value = "10.6.21.65"
$ws.Range("G39").Value = 300
$ws.Range("H39").Value = 45519.78981166667
$ws.Range("H39").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 40
$ws.Range("A40").Value = "Домик для Мышки Норушки"
$ws.Range("B40").Value = "без трубы"
$ws.Range("C40").Value = 2881
$ws.Range("D40").Value = 2885
$ws.Range("E40").Value = "10.6.21.66"
$ws.Range("F40").Value = "10.6.21.70"
$ws.Range("G40").Value = 5
$ws.Range("H40").Value = 45519.81906033565
$ws.Range("H40").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 41
$ws.Range("A41").Value = "Домик для Мышки Норушки"
$ws.Range("B41").Value = "без трубы"
$ws.Range("C41").Value = 2886
$ws.Range("D41").Value = 3145
$ws.Range("E41").Value = "10.6.21.71"
$ws.Range("F41").Value = "10.6.22.74"
$ws.Range("G41").Value = 260
$ws.Range("H41").Value = 45519.82235841767
$ws.Range("H41").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Rows added:" $ws.UsedRange.Rows.Count
